$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry updates one cell. AsText=$true is used for cells whose new value
# looks like a plain number (e.g. "242.46") so Excel does not silently convert
# the text into a numeric value and lose the exact original formatting/precision.
$updates = @(
    @{ Cell = "D2"; Value = "36.242.30"; AsText = $false },
    @{ Cell = "E2"; Value = "  -3.83%  "; AsText = $false },
    @{ Cell = "D3"; Value = "1.966.88"; AsText = $false },
    @{ Cell = "E3"; Value = "  -4.08%  "; AsText = $false },
    @{ Cell = "E4"; Value = "  +0.25%  "; AsText = $false },
    @{ Cell = "D5"; Value = "242.46"; AsText = $true },
    @{ Cell = "E5"; Value = "  -3.95%  "; AsText = $false },
    @{ Cell = "D6"; Value = "0.625"; AsText = $true },
    @{ Cell = "E6"; Value = "  -4.24%  "; AsText = $false },
    @{ Cell = "D7"; Value = "59.95"; AsText = $true },
    @{ Cell = "E7"; Value = "  -8.66%  "; AsText = $false },
    @{ Cell = "E8"; Value = "  +0.16%  "; AsText = $false },
    @{ Cell = "E9"; Value = "  -1.23%  "; AsText = $false },
    @{ Cell = "D10"; Value = "56.91"; AsText = $true },
    @{ Cell = "E10"; Value = "  -4.13%  "; AsText = $false },
    @{ Cell = "D11"; Value = "0.0800"; AsText = $true },
    @{ Cell = "E11"; Value = "  +5.25%  "; AsText = $false },
    @{ Cell = "E12"; Value = "  -0.53%  "; AsText = $false },
    @{ Cell = "D13"; Value = "0.859"; AsText = $true },
    @{ Cell = "E13"; Value = "  -6.83%  "; AsText = $false },
    @{ Cell = "D14"; Value = "22.38"; AsText = $true },
    @{ Cell = "E14"; Value = "  +9.28%  "; AsText = $false },
    @{ Cell = "D15"; Value = "14.06"; AsText = $true },
    @{ Cell = "E15"; Value = "  -7.51%  "; AsText = $false },
    @{ Cell = "D16"; Value = "2.255.95"; AsText = $false },
    @{ Cell = "E16"; Value = "  -3.99%  "; AsText = $false },
    @{ Cell = "E17"; Value = "  -2.90%  "; AsText = $false },
    @{ Cell = "D18"; Value = "1.965.04"; AsText = $false },
    @{ Cell = "E18"; Value = "  -4.33%  "; AsText = $false },
    @{ Cell = "D19"; Value = "36.131.37"; AsText = $false },
    @{ Cell = "E19"; Value = "  -3.78%  "; AsText = $false },
    @{ Cell = "E20"; Value = "  -3.70%  "; AsText = $false },
    @{ Cell = "E21"; Value = "  -2.27%  "; AsText = $false },
    @{ Cell = "D22"; Value = "236.85"; AsText = $true },
    @{ Cell = "E22"; Value = "  -0.57%  "; AsText = $false },
    @{ Cell = "E23"; Value = "  -2.78%  "; AsText = $false },
    @{ Cell = "E24"; Value = "  -0.05%  "; AsText = $false },
    @{ Cell = "D25"; Value = "2.54"; AsText = $true },
    @{ Cell = "E25"; Value = "  -5.28%  "; AsText = $false },
    @{ Cell = "D26"; Value = "2.28"; AsText = $true },
    @{ Cell = "E26"; Value = "  -4.39%  "; AsText = $false },
    @{ Cell = "D27"; Value = "9.81"; AsText = $true },
    @{ Cell = "E27"; Value = "  +1.74%  "; AsText = $false },
    @{ Cell = "D28"; Value = "160.69"; AsText = $true },
    @{ Cell = "E28"; Value = "  +0.18%  "; AsText = $false },
    @{ Cell = "D29"; Value = "19.82"; AsText = $true },
    @{ Cell = "E29"; Value = "  -0.78%  "; AsText = $false },
    @{ Cell = "D30"; Value = "0.127"; AsText = $true },
    @{ Cell = "E30"; Value = "  +11.73%  "; AsText = $false },
    @{ Cell = "E31"; Value = "  -2.35%  "; AsText = $false },
    @{ Cell = "E32"; Value = "  -7.72%  "; AsText = $false },
    @{ Cell = "E33"; Value = "  -5.82%  "; AsText = $false },
    @{ Cell = "E34"; Value = "  +0.67%  "; AsText = $false },
    @{ Cell = "E35"; Value = "  -7.37%  "; AsText = $false },
    @{ Cell = "D36"; Value = "6.28"; AsText = $true },
    @{ Cell = "E36"; Value = "  +5.74%  "; AsText = $false },
    @{ Cell = "E37"; Value = "  -6.40%  "; AsText = $false },
    @{ Cell = "E38"; Value = "  +0.32%  "; AsText = $false },
    @{ Cell = "E39"; Value = "  -1.76%  "; AsText = $false },
    @{ Cell = "D40"; Value = "3.09"; AsText = $true },
    @{ Cell = "E40"; Value = "  +9.19%  "; AsText = $false },
    @{ Cell = "D41"; Value = "0.0986"; AsText = $true },
    @{ Cell = "E41"; Value = "  -3.61%  "; AsText = $false },
    @{ Cell = "E42"; Value = "  -0.60%  "; AsText = $false },
    @{ Cell = "D43"; Value = "2.85"; AsText = $true },
    @{ Cell = "E43"; Value = "  -2.86%  "; AsText = $false },
    @{ Cell = "E44"; Value = "  -2.84%  "; AsText = $false },
    @{ Cell = "E45"; Value = "  -4.57%  "; AsText = $false },
    @{ Cell = "D46"; Value = "92.40"; AsText = $true },
    @{ Cell = "E46"; Value = "  -3.37%  "; AsText = $false },
    @{ Cell = "E47"; Value = "  -5.97%  "; AsText = $false },
    @{ Cell = "D48"; Value = "7.51"; AsText = $true },
    @{ Cell = "E48"; Value = "  -7.10%  "; AsText = $false },
    @{ Cell = "D49"; Value = "1.335.58"; AsText = $false },
    @{ Cell = "E49"; Value = "  -6.44%  "; AsText = $false },
    @{ Cell = "E50"; Value = "  -3.59%  "; AsText = $false },
    @{ Cell = "D51"; Value = "2.149.81"; AsText = $false },
    @{ Cell = "E51"; Value = "  -3.79%  "; AsText = $false }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    if ($u.AsText) {
        $cell.NumberFormat = "@"
        $cell.Value = $u.Value
        $cell.ClearFormats()
    } else {
        $cell.Value = $u.Value
    }
}
